# Regenerate Excel, and add back to list of files we check

$wb = $excel.ActiveWorkbook

# 1. Rename the main sheet from "Export this as TSV" to "Export as TSV"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Export as TSV"

# 2. Freeze the header row (row 1) on the main sheet
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Add errorTitle / error messages to all the data validations on the main sheet
$validations = @(
    @{ Sqref = "I2:I1048576"; Title = "Value must come from list"; Message = "Value must be one of: mass_spectrometry." },
    @{ Sqref = "J2:J1048576"; Title = "Value must come from list"; Message = "Value must be one of: LC-MS (metabolomics) / LC-MS/MS (label-free proteomics) / MS (shotgun lipidomics)." },
    @{ Sqref = "K2:K1048576"; Title = "Value must come from list"; Message = "Value must be one of: protein / metabolites / lipids." },
    @{ Sqref = "L2:L1048576"; Title = "Not a boolean"; Message = 'The values in this column must be "TRUE" or "FALSE".' },
    @{ Sqref = "Q2:Q1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Sqref = "R2:R1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Sqref = "AB2:AB1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Sqref = "AC2:AC1048576"; Title = "Value must come from list"; Message = "Value must be one of: um / mm / cm." },
    @{ Sqref = "AD2:AD1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Sqref = "AE2:AE1048576"; Title = "Value must come from list"; Message = "Value must be one of: C." },
    @{ Sqref = "AF2:AF1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Sqref = "AG2:AG1048576"; Title = "Value must come from list"; Message = "Value must be one of: um / mm / cm." },
    @{ Sqref = "AH2:AH1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Sqref = "AI2:AI1048576"; Title = "Value must come from list"; Message = "Value must be one of: nL/min / mL/min." }
)

foreach ($item in $validations) {
    $range = $ws.Range($item.Sqref)
    $dv = $range.Validation
    $dv.ErrorTitle = $item.Title
    $dv.ErrorMessage = $item.Message
}
